# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.865.30'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '1.802.48'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range("D5").Formula = '="309.26"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').Value = '  +4.17%  '
$ws.Range("D8").Formula = '="0.3701"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range('E8').Value = '  -2.04%  '
$ws.Range("D9").Formula = '="0.07385"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range("D10").Formula = '="0.8705"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range("D11").Formula = '="20.37"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range('E11').Value = '  -3.03%  '
$ws.Range('D12').Value = '1.821.37'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range("D13").Formula = '="5.360"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range("D14").Formula = '="92.41"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range("D15").Formula = '="6.482"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range('E15').Value = '  -3.87%  '
$ws.Range("D16").Formula = '="0.07023"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range("D18").Formula = '="0.000008709"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range("D19").Formula = '="1.000"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range("D20").Formula = '="14.70"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range('E20').Value = '  -3.04%  '
$ws.Range('D21').Value = '26.856.05'
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range("D22").Formula = '="5.294"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range("D23").Formula = '="10.63"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range('E23').Value = '  -3.32%  '
$ws.Range('D24').Value = '1.985.39'
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range("D25").Formula = '="1.898"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range('E25').Value = '  -3.70%  '
$ws.Range("D26").Formula = '="151.64"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range("D27").Formula = '="18.31"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range('E27').Value = '  -1.82%  '
$ws.Range("D28").Formula = '="2.137"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range('E28').Value = '  -8.16%  '
$ws.Range("D29").Formula = '="5.267"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range('E29').Value = '  -2.28%  '
$ws.Range("D30").Formula = '="115.74"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range('E30').Value = '  -1.84%  '
$ws.Range("D31").Formula = '="0.08946"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range("D32").Formula = '="0.7595"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range('E32').Value = '  -4.70%  '
$ws.Range("D33").Formula = '="1.151"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range('E33').Value = '  -4.25%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Formula = '="4.457"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Formula = '="2.913"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range("D36").Formula = '="0.9998"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range("D37").Formula = '="1.100"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range("D38").Formula = '="0.01957"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range('E38').Value = '  -1.49%  '
$ws.Range("D40").Formula = '="2.927"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range("D41").Formula = '="7.248"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range('E41').Value = '  -0.71%  '
$ws.Range('E42').Value = '  +2.40%  '
$ws.Range("D43").Formula = '="0.5287"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range("D44").Formula = '="0.1660"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range('E44').Value = '  -3.36%  '
$ws.Range("D45").Formula = '="8.495"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range('E45').Value = '  -2.09%  '
$ws.Range("D46").Formula = '="0.5002"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('E47').Value = '  -3.44%  '
$ws.Range("D48").Formula = '="104.03"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range("D49").Formula = '="0.9997"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range("D50").Formula = '="1.664"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range('E50').Value = '  -2.06%  '
$ws.Range("D51").Formula = '="0.06291"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range('E51').Value = '  -1.87%  '

$excel.CutCopyMode = 0

